$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4a"
$ws.Range("C2").Value = "Plxnb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 26.83824466666667
$ws.Range("H2").Value = 80.514734
$ws.Range("I2").Value = 0.8882651037973995
$ws.Range("J2").Value = 0.8882651037973996
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.448495
$ws.Range("N2").Value = 1.345485
$ws.Range("O2").Value = 0.05033609585745587
$ws.Range("P2").Value = 0.05033609585745586
$ws.Range("Q2").Value = 12.03681854177667
$ws.Range("R2").Value = 108.33136687599
$ws.Range("S2").Value = 0.04471179741157889
$ws.Range("T2").Value = 0.04471179741157889

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4a"
$ws.Range("C3").Value = "Plxnb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 26.83824466666667
$ws.Range("H3").Value = 80.514734
$ws.Range("I3").Value = 0.8882651037973995
$ws.Range("J3").Value = 0.8882651037973996
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6116653333333333
$ws.Range("N3").Value = 1.834996
$ws.Range("O3").Value = 0.06864924882406573
$ws.Range("P3").Value = 0.06864924882406571
$ws.Range("Q3").Value = 16.41602387011822
$ws.Range("R3").Value = 147.744214831064
$ws.Range("S3").Value = 0.06097873213232225
$ws.Range("T3").Value = 0.06097873213232224

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4a"
$ws.Range("C4").Value = "Plxnb1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 26.83824466666667
$ws.Range("H4").Value = 80.514734
$ws.Range("I4").Value = 0.8882651037973995
$ws.Range("J4").Value = 0.8882651037973996
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.849847333333334
$ws.Range("N4").Value = 23.549542
$ws.Range("O4").Value = 0.8810146553184784
$ws.Range("P4").Value = 0.8810146553184784
$ws.Range("Q4").Value = 210.6761233279809
$ws.Range("R4").Value = 1896.085109951828
$ws.Range("S4").Value = 0.7825745742534984
$ws.Range("T4").Value = 0.7825745742534985

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema4a"
$ws.Range("C5").Value = "Plxnb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.239011
$ws.Range("H5").Value = 3.717033
$ws.Range("I5").Value = 0.04100753414354395
$ws.Range("J5").Value = 0.04100753414354396
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.448495
$ws.Range("N5").Value = 1.345485
$ws.Range("O5").Value = 0.05033609585745587
$ws.Range("P5").Value = 0.05033609585745586
$ws.Range("Q5").Value = 0.555690238445
$ws.Range("R5").Value = 5.001212146005
$ws.Range("S5").Value = 0.002064159169527323
$ws.Range("T5").Value = 0.002064159169527323

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4a"
$ws.Range("C6").Value = "Plxnb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.239011
$ws.Range("H6").Value = 3.717033
$ws.Range("I6").Value = 0.04100753414354395
$ws.Range("J6").Value = 0.04100753414354396
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6116653333333333
$ws.Range("N6").Value = 1.834996
$ws.Range("O6").Value = 0.06864924882406573
$ws.Range("P6").Value = 0.06864924882406571
$ws.Range("Q6").Value = 0.7578600763186665
$ws.Range("R6").Value = 6.820740686867999
$ws.Range("S6").Value = 0.00281513641508152
$ws.Range("T6").Value = 0.00281513641508152

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4a"
$ws.Range("C7").Value = "Plxnb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.239011
$ws.Range("H7").Value = 3.717033
$ws.Range("I7").Value = 0.04100753414354395
$ws.Range("J7").Value = 0.04100753414354396
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.849847333333334
$ws.Range("N7").Value = 23.549542
$ws.Range("O7").Value = 0.8810146553184784
$ws.Range("P7").Value = 0.8810146553184784
$ws.Range("Q7").Value = 9.726047194320666
$ws.Range("R7").Value = 87.53442474888601
$ws.Range("S7").Value = 0.03612823855893511
$ws.Range("T7").Value = 0.03612823855893512

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema4a"
$ws.Range("C8").Value = "Plxnb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.136972666666667
$ws.Range("H8").Value = 6.410918000000001
$ws.Range("I8").Value = 0.07072736205905639
$ws.Range("J8").Value = 0.0707273620590564
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.448495
$ws.Range("N8").Value = 1.345485
$ws.Range("O8").Value = 0.05033609585745587
$ws.Range("P8").Value = 0.05033609585745586
$ws.Range("Q8").Value = 0.9584215561366668
$ws.Range("R8").Value = 8.62579400523
$ws.Range("S8").Value = 0.00356013927634965
$ws.Range("T8").Value = 0.00356013927634965

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema4a"
$ws.Range("C9").Value = "Plxnb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.136972666666667
$ws.Range("H9").Value = 6.410918000000001
$ws.Range("I9").Value = 0.07072736205905639
$ws.Range("J9").Value = 0.0707273620590564
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6116653333333333
$ws.Range("N9").Value = 1.834996
$ws.Range("O9").Value = 0.06864924882406573
$ws.Range("P9").Value = 0.06864924882406571
$ws.Range("Q9").Value = 1.307112098480889
$ws.Range("R9").Value = 11.764008886328
$ws.Range("S9").Value = 0.004855380276661948
$ws.Range("T9").Value = 0.004855380276661948

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema4a"
$ws.Range("C10").Value = "Plxnb1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.136972666666667
$ws.Range("H10").Value = 6.410918000000001
$ws.Range("I10").Value = 0.07072736205905639
$ws.Range("J10").Value = 0.0707273620590564
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.849847333333334
$ws.Range("N10").Value = 23.549542
$ws.Range("O10").Value = 0.8810146553184784
$ws.Range("P10").Value = 0.8810146553184784
$ws.Range("Q10").Value = 16.77490918883956
$ws.Range("R10").Value = 150.974182699556
$ws.Range("S10").Value = 0.0623118425060448
$ws.Range("T10").Value = 0.0623118425060448

